# Inserts one new data row (row 30) into the "Hortaliza, Femacal de La Calera - Haba"
# sheet, shifting the existing rows 30-135 down to 31-136, and populates the new
# row with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 30..135 down by one, creating a blank row 30 (format is inherited
# from the surrounding rows, so the date column keeps its date number format).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new observation.
$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "Femacal de La Calera"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44690
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = 100112026
$ws.Range("G30").Value = "Haba"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 73
$ws.Range("K30").Value = 19000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 19521
$ws.Range("N30").Value = "$/malla 25 kilos"
$ws.Range("O30").Value = "Provincia de Limarí"
$ws.Range("P30").Value = 781
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
